$d = $word.ActiveDocument

# 1. "questions in our mind" -> "questions on our mind"
$d.Content.Find.Execute(
  "questions in our mind",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "questions on our mind", 2) | Out-Null

# 2. "...go to? Maybe I would want to mingle..." -> "...go to. Maybe I would want to mingle..."
$d.Content.Find.Execute(
  "want to go to? Maybe I would want to mingle",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "want to go to. Maybe I would want to mingle", 2) | Out-Null

# 3. "...dislikes. Or possibly find..." -> "...dislikes? Or possibly find..."
$d.Content.Find.Execute(
  "having matching likes or dislikes. Or possibly find",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "having matching likes or dislikes? Or possibly find", 2) | Out-Null

# 4. "So where should we look ... The large public space itself!"
#    -> "So where should we be digging ... The answer is now obvious in this era of
#        ubiquitous computing - the large public space itself!"
$d.Content.Find.Execute(
  "So where should we look for this kind of information. The large public space itself!",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "So where should we be digging for this kind of information. The answer is now obvious in this era of ubiquitous computing - the large public space itself!", 2) | Out-Null

# 5. "solves certain problems of discovering people" -> "solves certain problems pertaining to discovering people"
$d.Content.Find.Execute(
  "solves certain problems of discovering people",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "solves certain problems pertaining to discovering people", 2) | Out-Null

# 6. "a certain physical range, but" -> "a certain broader area, but"
$d.Content.Find.Execute(
  "bring you within a certain physical range, but",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "bring you within a certain broader area, but", 2) | Out-Null

# 7. insert new sentence about co-ordinates + change "useful suggestions" -> "meaningful answers"
#    and "This would include" -> "These would include"
$d.Content.Find.Execute(
  "it is not “smart” data. Using contextual data from Google Beacon and processing user data on the cloud, we can provide useful suggestions to a user entering a public space. This would include",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "it is not “smart” data. Its co-ordinates do not have meaning attached to it.  Using contextual data from Google Beacon and processing user data on the cloud, we can provide meaningful answers to a user entering a public space. These would include", 2) | Out-Null

# 8. "If not the above, making" -> "If not one of the above, making"
$d.Content.Find.Execute(
  "If not the above, making new friends",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "If not one of the above, making new friends", 2) | Out-Null

# 9. "making a physical public smart" -> "making a physical public space smart"
$d.Content.Find.Execute(
  "making a physical public smart and knowledgeable",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "making a physical public space smart and knowledgeable", 2) | Out-Null

# 10. "would often need about the space." -> "would often need to know about that space."
$d.Content.Find.Execute(
  "would often need about the space.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "would often need to know about that space.", 2) | Out-Null

# 11. Add new sentence after "...and what you can hope to experience from the space. "
$d.Content.Find.Execute(
  "and what you can hope to experience from the space. ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "and what you can hope to experience from the space. The ability to discover your contact and finding the right person at the right time is always a very valuable asset to have. ", 2) | Out-Null

# Move the "_GoBack" bookmark from its old spot (between "like-" and "minded") to its
# new spot, right after "...discover your contact" in the sentence we just inserted.
$bmRange = $d.Content
$bmRange.Find.Execute("discover your contact", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmTarget = $d.Range($bmRange.End, $bmRange.End)
$d.Bookmarks.Add("_GoBack", $bmTarget) | Out-Null

# Delete the stray empty paragraph that sits between "...approach." and "2. Value
# Preposition (Potential Impact)" -- merges them the way the target document does.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
  $ptext = $paras.Item($i).Range.Text
  if ($ptext -eq "`r" -or $ptext -eq "") {
    $nextText = ""
    if (($i + 1) -le $paras.Count) {
      $nextText = $paras.Item($i + 1).Range.Text
    }
    $prevText = ""
    if (($i - 1) -ge 1) {
      $prevText = $paras.Item($i - 1).Range.Text
    }
    if ($nextText -like "2. Value Preposition*" -and $prevText -like "*easy to approach.*") {
      $paras.Item($i).Range.Delete() | Out-Null
      break
    }
  }
}
